$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.038.38"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "3.334.34"
$ws.Range("E3").Value = "  +0.52%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.11%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "178.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.72%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.593"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.79%  "
$ws.Range("E9").Value = "  +3.80%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.584"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.15%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "47.95"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.51%  "
$ws.Range("E12").Value = "  +1.74%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "697.18"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.85%  "
$ws.Range("D14").Value = "3.875.21"
$ws.Range("E14").Value = "  +0.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "8.44"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.70%  "
$ws.Range("D16").Value = "68.067.75"
$ws.Range("E16").Value = "  -0.16%  "
$ws.Range("E17").Value = "  +1.05%  "
$ws.Range("D18").Value = "3.331.87"
$ws.Range("E18").Value = "  +0.27%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.15"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.55%  "
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  -1.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "100.81"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.10%  "
$ws.Range("E25").Value = "  +1.71%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "33.16"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.81%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.03"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "574.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.62%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.03"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.106"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.66%  "
$ws.Range("D34").Value = "3.733.63"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "57.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.83%  "
$ws.Range("E36").Value = "  -0.05%  "
$ws.Range("E37").Value = "  +0.32%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "35.23"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +8.51%  "
$ws.Range("E39").Value = "  +1.13%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.39%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "
$ws.Range("E42").Value = "  +1.79%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.335"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.30"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0410"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.51%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.63"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.07%  "
$ws.Range("E47").Value = "  +0.97%  "
$ws.Range("E48").Value = "  -0.32%  "
$ws.Range("E49").Value = "  -1.17%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "130.55"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.17%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.35%  "
